$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as text (matches original inlineStr layout)
$ws.Range("A1:A40").NumberFormat = "@"

$data = @(
    ,@("20170822", 1154.0)
    ,@("20170904", 3400.0)
    ,@("20170929", 2275.0)
    ,@("20171013", 2902.0)
    ,@("20171206", 1696.0)
    ,@("20171228", 1960.0)
    ,@("20180108", 1812.0)
    ,@("20180116", 2116.0)
    ,@("20180207", 1344.0)
    ,@("20180312", 1685.0)
    ,@("20180326", 1398.0)
    ,@("20180413", 1990.0)
    ,@("20180423", 1639.0)
    ,@("20180528", 2257.0)
    ,@("20180706", 1450.0)
    ,@("20180713", 1705.0)
    ,@("20181012", 1092.0)
    ,@("20181113", 1453.0)
    ,@("20181126", 1291.0)
    ,@("20181211", 1909.0)
    ,@("20181228", 1590.0)
    ,@("20190109", 2009.0)
    ,@("20190130", 1610.0)
    ,@("20190225", 1936.0)
    ,@("20190311", 1660.0)
    ,@("20190319", 2184.0)
    ,@("20190606", 1398.0)
    ,@("20190724", 1776.0)
    ,@("20190812", 1351.0)
    ,@("20191121", 3360.0)
    ,@("20191230", 1648.0)
    ,@("20200121", 1833.0)
    ,@("20200204", 1323.0)
    ,@("20200225", 1946.0)
    ,@("20200319", 1419.0)
    ,@("20200327", 1800.0)
    ,@("20200428", 1430.0)
    ,@("20200521", 2668.0)
    ,@("20200529", 1983.0)
    ,@("20200804", 10800.0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
